# Auto-generated edit script applying numeric corrections to the leve-profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28 (Leve Item ID 27772)
$ws.Cells.Item(28, 8).Value = 255.25   # H28: was 256.25
$ws.Cells.Item(28, 9).Value = 255.25   # I28: was 256.25
$ws.Cells.Item(28, 11).Value = 255.25   # K28: was 256.25
$ws.Cells.Item(28, 13).Value = 229.75   # M28: was 228.75

# Row 31 (Leve Item ID 4576)
$ws.Cells.Item(31, 8).Value = 996.5454999999999   # H31: was 1089.9
$ws.Cells.Item(31, 9).Value = 996.5454999999999   # I31: was 1089.9
$ws.Cells.Item(31, 11).Value = 2989.6365   # K31: was 3269.7
$ws.Cells.Item(31, 13).Value = -2759.6365   # M31: was -3039.7

# Row 58 (Leve Item ID 4606)
$ws.Cells.Item(58, 8).Value = 367.14285   # H58: was 364.2857
$ws.Cells.Item(58, 9).Value = 289   # I58: was 311.25
$ws.Cells.Item(58, 10).Value = 562.5   # J58: was 435
$ws.Cells.Item(58, 11).Value = 867   # K58: was 933.75
$ws.Cells.Item(58, 12).Value = 1687.5   # L58: was 1305
$ws.Cells.Item(58, 13).Value = -717   # M58: was -783.75
$ws.Cells.Item(58, 14).Value = -1987.5   # N58: was -1605

# Row 64 (Leve Item ID 5506)
$ws.Cells.Item(64, 8).Value = 4509.1   # H64: was 4699.5
$ws.Cells.Item(64, 9).Value = 3824.25   # I64: was 3865.6667
$ws.Cells.Item(64, 10).Value = 4965.6665   # J64: was 5199.8
$ws.Cells.Item(64, 11).Value = 3824.25   # K64: was 3865.6667
$ws.Cells.Item(64, 12).Value = 4965.6665   # L64: was 5199.8
$ws.Cells.Item(64, 13).Value = -3576.25   # M64: was -3617.6667
$ws.Cells.Item(64, 14).Value = -5461.6665   # N64: was -5695.8

# Row 67 (Leve Item ID 5506)
$ws.Cells.Item(67, 8).Value = 4509.1   # H67: was 4699.5
$ws.Cells.Item(67, 9).Value = 3824.25   # I67: was 3865.6667
$ws.Cells.Item(67, 10).Value = 4965.6665   # J67: was 5199.8
$ws.Cells.Item(67, 11).Value = 3824.25   # K67: was 3865.6667
$ws.Cells.Item(67, 12).Value = 4965.6665   # L67: was 5199.8
$ws.Cells.Item(67, 13).Value = -2966.25   # M67: was -3007.6667
$ws.Cells.Item(67, 14).Value = -6681.6665   # N67: was -6915.8

# Row 103 (Leve Item ID 19909)
$ws.Cells.Item(103, 8).Value = 199   # H103: was 198.5
$ws.Cells.Item(103, 10).Value = 199   # J103: was 198.5
$ws.Cells.Item(103, 12).Value = 597   # L103: was 595.5
$ws.Cells.Item(103, 14).Value = -1769   # N103: was -1767.5

# Row 113 (Leve Item ID 27775)
$ws.Cells.Item(113, 8).Value = 4356.143   # H113: was 4398.8
$ws.Cells.Item(113, 9).Value = 4298.6   # I113: was 4373.5
$ws.Cells.Item(113, 11).Value = 4298.6   # K113: was 4373.5
$ws.Cells.Item(113, 13).Value = -1044.6   # M113: was -1119.5

# Row 138 (Leve Item ID 44169)
$ws.Cells.Item(138, 8).Value = 3277   # H138: was 3404.5
$ws.Cells.Item(138, 9).Value = 2554.8   # I138: was 2733.2222
$ws.Cells.Item(138, 10).Value = 3999.2   # J138: was 3953.7273
$ws.Cells.Item(138, 11).Value = 7664.400000000001   # K138: was 8199.6666
$ws.Cells.Item(138, 12).Value = 11997.6   # L138: was 11861.1819
$ws.Cells.Item(138, 13).Value = -2524.400000000001   # M138: was -3059.6666
$ws.Cells.Item(138, 14).Value = -22277.6   # N138: was -22141.1819

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Cells.Item(32, 8).Value = 2335603   # H32: was 2266915.5
$ws.Cells.Item(32, 9).Value = 2189840.8   # I32: was 2123488.8
$ws.Cells.Item(32, 11).Value = 2189840.8   # K32: was 2123488.8
$ws.Cells.Item(32, 13).Value = -2189553.8   # M32: was -2123201.8

# Row 122 (Leve Item ID 36168)
$ws.Cells.Item(122, 8).Value = 3559.8   # H122: was 3649.5
$ws.Cells.Item(122, 10).Value = 3574.75   # J122: was 3799
$ws.Cells.Item(122, 12).Value = 10724.25   # L122: was 11397
$ws.Cells.Item(122, 14).Value = -15624.25   # N122: was -16297

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (Leve Item ID 12526)
$ws.Cells.Item(86, 8).Value = 1500   # H86: was 1876.25
$ws.Cells.Item(86, 9).Value = 1500   # I86: was 1666
$ws.Cells.Item(86, 10).Value = 0   # J86: was 2507
$ws.Cells.Item(86, 11).Value = 1500   # K86: was 1666
$ws.Cells.Item(86, 12).Value = 0   # L86: was 2507
$ws.Cells.Item(86, 14).Value = $null   # N86: was -4753

# Row 89 (Leve Item ID 12526)
$ws.Cells.Item(89, 8).Value = 1500   # H89: was 1876.25
$ws.Cells.Item(89, 9).Value = 1500   # I89: was 1666
$ws.Cells.Item(89, 10).Value = 0   # J89: was 2507
$ws.Cells.Item(89, 11).Value = 7500   # K89: was 8330
$ws.Cells.Item(89, 12).Value = 0   # L89: was 12535
$ws.Cells.Item(89, 14).Value = $null   # N89: was -23767

# Row 94 (Leve Item ID 19939)
$ws.Cells.Item(94, 8).Value = 532.2857   # H94: was 398.6
$ws.Cells.Item(94, 9).Value = 537.6667   # I94: was 425.125
$ws.Cells.Item(94, 10).Value = 500   # J94: was 292.5
$ws.Cells.Item(94, 11).Value = 537.6667   # K94: was 425.125
$ws.Cells.Item(94, 12).Value = 500   # L94: was 292.5
$ws.Cells.Item(94, 13).Value = -86.66669999999999   # M94: was 25.875
$ws.Cells.Item(94, 14).Value = -1402   # N94: was -1194.5

# Row 105 (Leve Item ID 19947)
$ws.Cells.Item(105, 8).Value = 2582.8333   # H105: was 2180.25
$ws.Cells.Item(105, 9).Value = 2499.4   # I105: was 2063.1428
$ws.Cells.Item(105, 11).Value = 2499.4   # K105: was 2063.1428
$ws.Cells.Item(105, 13).Value = -752.4000000000001   # M105: was -316.1428000000001

# Row 134 (Leve Item ID 43998)
$ws.Cells.Item(134, 8).Value = 2060.3333   # H134: was 2100.3572
$ws.Cells.Item(134, 9).Value = 1992.2142   # I134: was 2030.0769
$ws.Cells.Item(134, 11).Value = 5976.642599999999   # K134: was 6090.2307
$ws.Cells.Item(134, 13).Value = -3441.642599999999   # M134: was -3555.2307

$ws = $wb.Worksheets.Item("CRP")
# Row 4 (Leve Item ID 3742)
$ws.Cells.Item(4, 8).Value = 241   # H4: was 244.1
$ws.Cells.Item(4, 9).Value = 217.54546   # I4: was 224
$ws.Cells.Item(4, 10).Value = 499   # J4: was 324.5
$ws.Cells.Item(4, 11).Value = 217.54546   # K4: was 224
$ws.Cells.Item(4, 12).Value = 499   # L4: was 324.5
$ws.Cells.Item(4, 13).Value = -105.54546   # M4: was -112
$ws.Cells.Item(4, 14).Value = -723   # N4: was -548.5

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (Leve Item ID 4650)
$ws.Cells.Item(4, 8).Value = 12222376   # H4: was 11000143
$ws.Cells.Item(4, 10).Value = 0   # J4: was 43
$ws.Cells.Item(4, 12).Value = 0   # L4: was 129
$ws.Cells.Item(4, 14).Value = $null   # N4: was -353

# Row 47 (Leve Item ID 4663)
$ws.Cells.Item(47, 8).Value = 928.4286   # H47: was 270.66666
$ws.Cells.Item(47, 9).Value = 928.4286   # I47: was 224.8
$ws.Cells.Item(47, 10).Value = 0   # J47: was 500
$ws.Cells.Item(47, 11).Value = 2785.2858   # K47: was 674.4000000000001
$ws.Cells.Item(47, 12).Value = 0   # L47: was 1500
$ws.Cells.Item(47, 14).Value = $null   # N47: was -2362

# Row 75 (Leve Item ID 12863)
$ws.Cells.Item(75, 8).Value = 3980   # H75: was 2300
$ws.Cells.Item(75, 10).Value = 3300   # J75: was 950
$ws.Cells.Item(75, 12).Value = 9900   # L75: was 2850
$ws.Cells.Item(75, 14).Value = -11896   # N75: was -4846

# Row 78 (Leve Item ID 12863)
$ws.Cells.Item(78, 8).Value = 3980   # H78: was 2300
$ws.Cells.Item(78, 10).Value = 3300   # J78: was 950
$ws.Cells.Item(78, 12).Value = 29700   # L78: was 8550
$ws.Cells.Item(78, 14).Value = -39684   # N78: was -18534

# Row 119 (Leve Item ID 27873)
$ws.Cells.Item(119, 8).Value = 3499   # H119: was 3332.3333
$ws.Cells.Item(119, 9).Value = 3499   # I119: was 3332.3333
$ws.Cells.Item(119, 11).Value = 10497   # K119: was 9996.999899999999
$ws.Cells.Item(119, 13).Value = -5659   # M119: was -5158.999899999999

$ws = $wb.Worksheets.Item("GSM")
# Row 46 (Leve Item ID 2078)
$ws.Cells.Item(46, 8).Value = 9311.111000000001   # H46: was 8108.3335
$ws.Cells.Item(46, 9).Value = 3800   # I46: was 2900
$ws.Cells.Item(46, 10).Value = 10000   # J46: was 9150
$ws.Cells.Item(46, 11).Value = 3800   # K46: was 2900
$ws.Cells.Item(46, 12).Value = 10000   # L46: was 9150
$ws.Cells.Item(46, 13).Value = -3644   # M46: was -2744
$ws.Cells.Item(46, 14).Value = -10312   # N46: was -9462

# Row 113 (Leve Item ID 27710)
$ws.Cells.Item(113, 8).Value = 724   # H113: was 699
$ws.Cells.Item(113, 9).Value = 449   # I113: was 499
$ws.Cells.Item(113, 11).Value = 449   # K113: was 499
$ws.Cells.Item(113, 13).Value = 1721   # M113: was 1671

# Row 131 (Leve Item ID 34747)
$ws.Cells.Item(131, 8).Value = 5000   # H131: was 0
$ws.Cells.Item(131, 9).Value = 5000   # I131: was 0
$ws.Cells.Item(131, 11).Value = 5000   # K131: was 0
$ws.Cells.Item(131, 13).Value = 40   # M131: newly added

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Cells.Item(7, 8).Value = 7174.6   # H7: was 7849.1113
$ws.Cells.Item(7, 9).Value = 6785.5713   # I7: was 9058.200000000001
$ws.Cells.Item(7, 11).Value = 6785.5713   # K7: was 9058.200000000001
$ws.Cells.Item(7, 13).Value = -6673.5713   # M7: was -8946.200000000001

# Row 40 (Leve Item ID 36248)
$ws.Cells.Item(40, 8).Value = 4457.1665   # H40: was 4750.25
$ws.Cells.Item(40, 9).Value = 3247.6667   # I40: was 4200.4
$ws.Cells.Item(40, 11).Value = 3247.6667   # K40: was 4200.4
$ws.Cells.Item(40, 13).Value = -3111.6667   # M40: was -4064.4

# Row 46 (Leve Item ID 5282)
$ws.Cells.Item(46, 8).Value = 1804.8889   # H46: was 2224.4
$ws.Cells.Item(46, 10).Value = 1982.5   # J46: was 2556.4285
$ws.Cells.Item(46, 12).Value = 1982.5   # L46: was 2556.4285
$ws.Cells.Item(46, 14).Value = -2358.5   # N46: was -2932.4285

# Row 61 (Leve Item ID 27740)
$ws.Cells.Item(61, 8).Value = 1416.4445   # H61: was 1542.7142
$ws.Cells.Item(61, 9).Value = 1399.8572   # I61: was 1466.6666
$ws.Cells.Item(61, 10).Value = 1474.5   # J61: was 1999
$ws.Cells.Item(61, 11).Value = 1399.8572   # K61: was 1466.6666
$ws.Cells.Item(61, 12).Value = 1474.5   # L61: was 1999
$ws.Cells.Item(61, 13).Value = -1197.8572   # M61: was -1264.6666
$ws.Cells.Item(61, 14).Value = -1878.5   # N61: was -2403

# Row 82 (Leve Item ID 12565)
$ws.Cells.Item(82, 8).Value = 916.1667   # H82: was 1049.5
$ws.Cells.Item(82, 10).Value = 1049.5   # J82: was 1249.5
$ws.Cells.Item(82, 12).Value = 1049.5   # L82: was 1249.5
$ws.Cells.Item(82, 14).Value = -1771.5   # N82: was -1971.5

# Row 85 (Leve Item ID 12565)
$ws.Cells.Item(85, 8).Value = 916.1667   # H85: was 1049.5
$ws.Cells.Item(85, 10).Value = 1049.5   # J85: was 1249.5
$ws.Cells.Item(85, 12).Value = 1049.5   # L85: was 1249.5
$ws.Cells.Item(85, 14).Value = -3545.5   # N85: was -3745.5

# Row 103 (Leve Item ID 18526)
$ws.Cells.Item(103, 8).Value = 17701   # H103: was 17901
$ws.Cells.Item(103, 10).Value = 17701   # J103: was 17901
$ws.Cells.Item(103, 12).Value = 17701   # L103: was 17901
$ws.Cells.Item(103, 14).Value = -20045   # N103: was -20245

# Row 113 (Leve Item ID 27740)
$ws.Cells.Item(113, 8).Value = 1416.4445   # H113: was 1542.7142
$ws.Cells.Item(113, 9).Value = 1399.8572   # I113: was 1466.6666
$ws.Cells.Item(113, 10).Value = 1474.5   # J113: was 1999
$ws.Cells.Item(113, 11).Value = 1399.8572   # K113: was 1466.6666
$ws.Cells.Item(113, 12).Value = 1474.5   # L113: was 1999
$ws.Cells.Item(113, 13).Value = 770.1428000000001   # M113: was 703.3334
$ws.Cells.Item(113, 14).Value = -5814.5   # N113: was -6339

# Row 126 (Leve Item ID 36249)
$ws.Cells.Item(126, 8).Value = 7174.6   # H126: was 7849.1113
$ws.Cells.Item(126, 9).Value = 6785.5713   # I126: was 9058.200000000001
$ws.Cells.Item(126, 11).Value = 20356.7139   # K126: was 27174.6
$ws.Cells.Item(126, 13).Value = -17886.7139   # M126: was -24704.6

# Row 132 (Leve Item ID 44058)
$ws.Cells.Item(132, 8).Value = 3165.4285   # H132: was 1643.2858
$ws.Cells.Item(132, 9).Value = 2031.6   # I132: was 1799.6666
$ws.Cells.Item(132, 10).Value = 6000   # J132: was 705
$ws.Cells.Item(132, 11).Value = 6094.799999999999   # K132: was 5398.9998
$ws.Cells.Item(132, 12).Value = 18000   # L132: was 2115
$ws.Cells.Item(132, 13).Value = -3564.799999999999   # M132: was -2868.9998
$ws.Cells.Item(132, 14).Value = -23060   # N132: was -7175

$ws = $wb.Worksheets.Item("WVR")
# Row 2 (Leve Item ID 3307)
$ws.Cells.Item(2, 8).Value = 314332.34   # H2: was 270499.16
$ws.Cells.Item(2, 10).Value = 314332.34   # J2: was 270499.16
$ws.Cells.Item(2, 12).Value = 314332.34   # L2: was 270499.16
$ws.Cells.Item(2, 14).Value = -314556.34   # N2: was -270723.16

# Row 27 (Leve Item ID 27174)
$ws.Cells.Item(27, 8).Value = 39990   # H27: was 39989.5
$ws.Cells.Item(27, 10).Value = 39990   # J27: was 39989.5
$ws.Cells.Item(27, 12).Value = 39990   # L27: was 39989.5
$ws.Cells.Item(27, 14).Value = -40128   # N27: was -40127.5

# Row 107 (Leve Item ID 27746)
$ws.Cells.Item(107, 8).Value = 3998.6365   # H107: was 4047.9167
$ws.Cells.Item(107, 9).Value = 3923.375   # I107: was 3997.4443
$ws.Cells.Item(107, 11).Value = 11770.125   # K107: was 11992.3329
$ws.Cells.Item(107, 13).Value = -9850.125   # M107: was -10072.3329

# Row 136 (Leve Item ID 44031)
$ws.Cells.Item(136, 8).Value = 1965   # H136: was 1965.625
$ws.Cells.Item(136, 9).Value = 2029.3334   # I136: was 2030
$ws.Cells.Item(136, 11).Value = 6088.0002   # K136: was 6090
$ws.Cells.Item(136, 13).Value = -3538.0002   # M136: was -3540

Write-Output "Applied leve-profit updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets."